$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3846153846153846
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 0.4347826086956522

$ws.Range("B3").Value = 0.5454545454545454
$ws.Range("C3").Value = 0.4285714285714285
$ws.Range("D3").Value = 0.4799999999999999

$ws.Range("B4").Value = 0.4583333333333333
$ws.Range("C4").Value = 0.4583333333333333
$ws.Range("D4").Value = 0.4583333333333333
$ws.Range("E4").Value = 0.4583333333333333

$ws.Range("B5").Value = 0.465034965034965
$ws.Range("C5").Value = 0.4642857142857143
$ws.Range("D5").Value = 0.457391304347826

$ws.Range("B6").Value = 0.4784382284382284
$ws.Range("C6").Value = 0.4583333333333333
$ws.Range("D6").Value = 0.4611594202898551

$ws.Range("B7").Value = 0.5
$ws.Range("D7").Value = 0.5

$ws.Range("B8").Value = 0.6428571428571429
$ws.Range("C8").Value = 0.6428571428571429
$ws.Range("D8").Value = 0.6428571428571429

$ws.Range("B9").Value = 0.5833333333333334
$ws.Range("C9").Value = 0.5833333333333334
$ws.Range("D9").Value = 0.5833333333333334
$ws.Range("E9").Value = 0.5833333333333334

$ws.Range("B10").Value = 0.5714285714285714
$ws.Range("C10").Value = 0.5714285714285714
$ws.Range("D10").Value = 0.5714285714285714

$ws.Range("B11").Value = 0.5833333333333334
$ws.Range("C11").Value = 0.5833333333333334
$ws.Range("D11").Value = 0.5833333333333334

$ws.Range("B12").Value = 0.6
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.6

$ws.Range("B13").Value = 0.7142857142857143
$ws.Range("C13").Value = 0.7142857142857143
$ws.Range("D13").Value = 0.7142857142857143

$ws.Range("B15").Value = 0.6571428571428571
$ws.Range("C15").Value = 0.6571428571428571
$ws.Range("D15").Value = 0.6571428571428571

$ws.Range("B16").Value = 0.6666666666666666
$ws.Range("D16").Value = 0.6666666666666666

$ws.Range("B17").Value = 0.4166666666666667
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 0.4545454545454545

$ws.Range("B18").Value = 0.5833333333333334
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 0.5384615384615384

$ws.Range("B19").Value = 0.5
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = 0.5
$ws.Range("E19").Value = 0.5

$ws.Range("B20").Value = 0.5
$ws.Range("C20").Value = 0.5
$ws.Range("D20").Value = 0.4965034965034965

$ws.Range("B21").Value = 0.513888888888889
$ws.Range("C21").Value = 0.5
$ws.Range("D21").Value = 0.5034965034965034

$ws.Range("B22").Value = 0.4615384615384616
$ws.Range("C22").Value = 0.6
$ws.Range("D22").Value = 0.5217391304347826

$ws.Range("B23").Value = 0.6363636363636364
$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = 0.5600000000000001

$ws.Range("B24").Value = 0.5416666666666666
$ws.Range("C24").Value = 0.5416666666666666
$ws.Range("D24").Value = 0.5416666666666666
$ws.Range("E24").Value = 0.5416666666666666

$ws.Range("B25").Value = 0.548951048951049
$ws.Range("C25").Value = 0.55
$ws.Range("D25").Value = 0.5408695652173914

$ws.Range("B26").Value = 0.5635198135198135
$ws.Range("C26").Value = 0.5416666666666666
$ws.Range("D26").Value = 0.5440579710144928
